# Applies the "Updating numbers in the PPT" commit:
#   - Slide 4 ("Then let's see Orchard Core!"): commit/issue counters
#     bumped from 6900/6100 to 7000/6200 (i.e. the leading "69"/"61"
#     runs become "70"/"62").
#   - The auto-updating "Update automatically" date footer field
#     (cached text "8/22/2023") is refreshed to "9/27/2023" everywhere
#     it is defined: the slide master, all 11 slide layouts, and the
#     notes master.

$p = $ppt.ActivePresentation

# ---------------------------------------------------------------
# 1. Slide 4 stats: 6900 -> 7000 commits, 6100 -> 6200 issues
# ---------------------------------------------------------------
$s4 = $p.Slides.Item(4)

for ($i = 1; $i -le $s4.Shapes.Count; $i++) {
    $shp = $s4.Shapes.Item($i)
    if (-not $shp.HasTextFrame) { continue }
    $tr = $shp.TextFrame.TextRange

    $t = $tr.Text
    if ($t.IndexOf("6900") -ge 0) {
        $idx = $t.IndexOf("6900")
        $tr.Characters($idx + 1, 2).Text = "70"
    }

    $t = $tr.Text
    if ($t.IndexOf("6100") -ge 0) {
        $idx = $t.IndexOf("6100")
        $tr.Characters($idx + 1, 2).Text = "62"
    }
}

# ---------------------------------------------------------------
# 2. Refresh the cached "datetimeFigureOut" footer text wherever a
#    "Date Placeholder" shape exists: slide master, every custom
#    layout, and the notes master.
# ---------------------------------------------------------------
function Update-DatePlaceholder($shapes, $newDate) {
    for ($j = 1; $j -le $shapes.Count; $j++) {
        $shp = $shapes.Item($j)
        if ($shp.Name -like "Date Placeholder*" -and $shp.HasTextFrame) {
            $tr = $shp.TextFrame.TextRange
            if ($tr.Length -gt 0) {
                $tr.Characters(1, $tr.Length).Text = $newDate
            } else {
                $tr.Text = $newDate
            }
        }
    }
}

$newDate = "9/27/2023"

$master = $p.SlideMaster
Update-DatePlaceholder $master.Shapes $newDate

$layouts = $master.CustomLayouts
for ($k = 1; $k -le $layouts.Count; $k++) {
    $layout = $layouts.Item($k)
    Update-DatePlaceholder $layout.Shapes $newDate
}

$notesMaster = $p.NotesMaster
Update-DatePlaceholder $notesMaster.Shapes $newDate
